$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts old E..H to F..I)
$ws.Columns.Item(5).Insert()

# New header for column E (copy header formatting: bold + centered)
$ws.Cells.Item(1,5).Value = "variacao_2024_2023_corrigida"
$ws.Cells.Item(1,5).Font.Bold = $true
$ws.Cells.Item(1,5).HorizontalAlignment = -4108

# New header for column J (margem_despesa), appended after column I
$ws.Cells.Item(1,10).Value = "margem_despesa"
$ws.Cells.Item(1,10).Font.Bold = $true
$ws.Cells.Item(1,10).HorizontalAlignment = -4108

# Per-row data: column E (variacao_2024_2023_corrigida) and column J (margem_despesa)
$ws.Cells.Item(2,5).Value = 0.07217724455430408
$ws.Cells.Item(2,10).Value = 0.09882407118801285
$ws.Cells.Item(3,5).Value = 0.05687339318765683
$ws.Cells.Item(3,10).Value = 0.08811137523135978
$ws.Cells.Item(4,5).Value = 0.007747811020725637
$ws.Cells.Item(4,10).Value = 0.05372346771450795
$ws.Cells.Item(5,5).Value = 0.1099054456163597
$ws.Cells.Item(5,10).Value = 0.1252338119314518
$ws.Cells.Item(6,5).Value = -0.02474494835859653
$ws.Cells.Item(6,10).Value = 0.0483
$ws.Cells.Item(7,5).Value = 0.1412448032998892
$ws.Cells.Item(7,10).Value = 0.1471713623099225
$ws.Cells.Item(8,5).Value = -0.070301201882641
$ws.Cells.Item(8,10).Value = 0.0483
$ws.Cells.Item(9,5).Value = 0.02754410892825732
$ws.Cells.Item(9,10).Value = 0.06207205446412866
$ws.Cells.Item(10,5).Value = 0.003774324351723868
$ws.Cells.Item(10,10).Value = 0.05094202704620671
$ws.Cells.Item(11,5).Value = 0.1041603739215926
$ws.Cells.Item(11,10).Value = 0.1212122617451148
$ws.Cells.Item(12,5).Value = 0.0833872618131637
$ws.Cells.Item(12,10).Value = 0.1066710832692146
$ws.Cells.Item(13,5).Value = 0.05872842322164051
$ws.Cells.Item(13,10).Value = 0.08940989625514836
$ws.Cells.Item(14,5).Value = -0.008622525589663144
$ws.Cells.Item(14,10).Value = 0.0483
$ws.Cells.Item(15,5).Value = 0.04799685304602352
$ws.Cells.Item(15,10).Value = 0.08189779713221647
$ws.Cells.Item(16,5).Value = 0.06807493135104581
$ws.Cells.Item(16,10).Value = 0.09595245194573207
$ws.Cells.Item(17,5).Value = 0.07650302863077885
$ws.Cells.Item(17,10).Value = 0.1018521200415452
$ws.Cells.Item(18,5).Value = 0.04214346369921396
$ws.Cells.Item(18,10).Value = 0.06937173184960699
$ws.Cells.Item(19,5).Value = 0.2171278419413167
$ws.Cells.Item(19,10).Value = 0.2002894893589217
$ws.Cells.Item(20,5).Value = -0.0587844644861697
$ws.Cells.Item(20,10).Value = 0.0483
$ws.Cells.Item(21,5).Value = 0.120880003691471
$ws.Cells.Item(21,10).Value = 0.1329160025840297
$ws.Cells.Item(22,5).Value = 0.1185231531470266
$ws.Cells.Item(22,10).Value = 0.1312662072029186
$ws.Cells.Item(23,5).Value = 0.06235166936010317
$ws.Cells.Item(23,10).Value = 0.0794758346800516
$ws.Cells.Item(24,5).Value = 0.1009291003082993
$ws.Cells.Item(24,10).Value = 0.1189503702158095
$ws.Cells.Item(25,5).Value = 0.05282282983040654
$ws.Cells.Item(25,10).Value = 0.07471141491520328
$ws.Cells.Item(26,5).Value = 0.116571484716731
$ws.Cells.Item(26,10).Value = 0.1299000393017117
$ws.Cells.Item(27,5).Value = 0.09844586889032292
$ws.Cells.Item(27,10).Value = 0.117212108223226
$ws.Cells.Item(28,5).Value = 0.07476567041636573
$ws.Cells.Item(28,10).Value = 0.08568283520818287

# Minor floating point recalculation updates to column D (total_receita_2023_corrigida)
$ws.Cells.Item(10,4).Value = 79598557003.43983
$ws.Cells.Item(12,4).Value = 7751748405.048451
$ws.Cells.Item(27,4).Value = 36466069939.88294

# Ensure dimension reflects full used range
Write-Output $ws.UsedRange.Address()
